# Regenerate merged AHB files
# 1. Rename header labels from _old/_new suffixes to _FV2410/_FV2504
# 2. Turn the used range into an Excel Table ("Table1")
# 3. Freeze the header row (pane split under row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels -------------------------------------------------
# "_old" columns (A:J) become "_FV2410"; "_new" columns (L:U) become "_FV2504".
# Column K ("diff") is left untouched.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    if ($text -ne $null) {
        if ($text.EndsWith("_old")) {
            $cell.Value2 = $text.Substring(0, $text.Length - 4) + "_FV2410"
        } elseif ($text.EndsWith("_new")) {
            $cell.Value2 = $text.Substring(0, $text.Length - 4) + "_FV2504"
        }
    }
}

# --- 2. Create a table over the used range ----------------------------------------
$range = $ws.Range("A1:U93")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
